$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clean slate for the whole table so the old "date" number format doesn't
# get dragged along while we re-style things below.
$ws.Range("A1:E4").ClearFormats()

# Center every cell in the table (header row + all data rows).
$ws.Range("A1:E4").HorizontalAlignment = -4108

# The pickup/drop-off date columns now hold plain text, not real dates.
$ws.Range("B2:B4").NumberFormat = "@"
$ws.Range("D2:D4").NumberFormat = "@"

# New pickup/drop-off date values, entered as plain text (not real dates).
$ws.Range("B2").Value = "07/20/2020"
$ws.Range("B3").Value = "07/25/2020"
$ws.Range("B4").Value = "07/22/2020"

$ws.Range("D2").Value = "08/30/2020"
$ws.Range("D3").Value = "08/30/2020"
$ws.Range("D4").Value = "07/30/2020"

# Widen columns to fit the new data, and make room for a new column F.
$ws.Columns("B").ColumnWidth = 14.28
$ws.Columns("C").ColumnWidth = 15.5
$ws.Columns("F").ColumnWidth = 11.28

# Move the selection, matching where the author clicked next.
$ws.Range("C8").Select()
